$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.802.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "'2.315.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'232.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'0.645"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("D7").Value = "'65.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.79%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.442"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.0972"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").Value = "'57.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "'26.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.08%  "
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "'2.638.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "'15.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "'6.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "'0.834"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "'2.299.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").Value = "'43.736.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "'0.0₃0977"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").Value = "'73.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").Value = "'6.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").Value = "'249.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'3.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.60%  "
$ws.Range("D26").Value = "'2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'2.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").Value = "'9.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").Value = "'174.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("D30").Value = "'22.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.26%  "
$ws.Range("D31").Value = "'0.132"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("D32").Value = "'1.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.26%  "
$ws.Range("D33").Value = "'0.126"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").Value = "'4.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.48%  "
$ws.Range("D35").Value = "'0.0685"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "'4.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.10%  "
$ws.Range("D37").Value = "'6.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'3.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D39").Value = "'2.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.36%  "
$ws.Range("D40").Value = "'0.0252"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'8.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.00%  "
$ws.Range("D43").Value = "'17.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "'4.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "'98.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.0950"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").Value = "'10.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.81%  "
$ws.Range("D48").Value = "'1.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "'0.000208"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'1.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.55%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.441.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
